$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing the existing row 3 (MuSCs data) down to row 4.
$ws.Rows("3").Insert()

# Row 2 (Sending=MuSCs, Ligand=Gdnf, Receptor=Gfra2, Target=ECs) - new TPM-based values
$ws.Range("D2").Value = "ECs"
$ws.Range("M2").Value = 0.9341033333333334
$ws.Range("N2").Value = 2.80231
$ws.Range("O2").Value = 0.04508188809474321
$ws.Range("P2").Value = 0.04508188809474322
$ws.Range("Q2").Value = 0.65192659609
$ws.Range("R2").Value = 5.867339364810001
$ws.Range("S2").Value = 0.04508188809474321
$ws.Range("T2").Value = 0.04508188809474322

# Row 3 (new, Target=FAPs) - carries the values that used to be in row 2
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Gfra2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.697917
$ws.Range("H3").Value = 2.093751
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.837667
$ws.Range("N3").Value = 32.513001
$ws.Range("O3").Value = 0.5230497242297513
$ws.Range("P3").Value = 0.5230497242297513
$ws.Range("Q3").Value = 7.563792039639001
$ws.Range("R3").Value = 68.074128356751
$ws.Range("S3").Value = 0.5230497242297513
$ws.Range("T3").Value = 0.5230497242297513

# Row 4 (Target=MuSCs) - data unchanged, but specificity columns recalculated
$ws.Range("O4").Value = 0.4318683876755055
$ws.Range("P4").Value = 0.4318683876755055
$ws.Range("S4").Value = 0.4318683876755055
$ws.Range("T4").Value = 0.4318683876755055
